$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q3" worksheet right before the existing
#    "2022-Q2" sheet (mirrors the workbook.xml sheet-order change in the
#    diff: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3).
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Match the look of the sibling fund-holding sheets: bold, centered,
# top-aligned, thin-bordered header row (style used on B1:H1 in every other
# quarter sheet) and the same styling on the numeric index column (A).
$hdr = $q3.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$idx = $q3.Range("A2:A3")
$idx.Font.Bold = $true
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160
$idx.Borders.LineStyle = 1

$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Fund figures in this table are stored as text, not numbers, matching the
# sibling quarter sheets (e.g. "16.84" as a string) - force text format
# before writing so the numeric-looking strings are not auto-converted.
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "517160"
$q3.Cells.Item(2,3).Value = "南方中证长江保护主题ETF"
$q3.Cells.Item(2,4).Value = "16.84"
$q3.Cells.Item(2,5).Value = "99.30"
$q3.Cells.Item(2,6).Value = "2.34"
$q3.Cells.Item(2,7).Value = "0.3941"
$q3.Cells.Item(2,8).Value = 8

$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "517330"
$q3.Cells.Item(3,3).Value = "易方达中证长江保护主题ETF"
$q3.Cells.Item(3,4).Value = "16.37"
$q3.Cells.Item(3,5).Value = "99.46"
$q3.Cells.Item(3,6).Value = "2.33"
$q3.Cells.Item(3,7).Value = "0.3814"
$q3.Cells.Item(3,8).Value = 8

# ---------------------------------------------------------------------------
# 2. Add the new quarter's summary row to the "总计" sheet: insert a row at
#    row 2 (pushing 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 down by one) and
#    fill it in with the 2022-Q3 totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Match formatting of the other data rows (column A carries the numbered
# index style; B:D stay unstyled like their neighbours).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.78

# Renumber the index column (A) sequentially for the rows that shifted down.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection (the workbook was showing
# "2021-Q3" before this edit).
# ---------------------------------------------------------------------------
$null = $wb.Worksheets.Item("2021-Q3").Activate()
$null = $wb.Worksheets.Item("2021-Q3").Range("A1").Select()
